$p = $ppt.ActivePresentation
$dsgs = $p.Designs
Write-Output ("Designs.Count=" + $dsgs.Count)
$d1 = $dsgs.Item(1)
Write-Output ("Design1 name=" + $d1.Name)
$sm = $d1.SlideMaster
Write-Output ("SlideMaster name=" + $sm.Name)
$th = $sm.Theme
Write-Output ("Theme=" + $th)
$tv = $th.ThemeVariants
Write-Output ("ThemeVariants count=" + $tv.Count)
